$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Execute Code" block (K5:L5 header) gains a new shortcut row.
# "debug cell" / "ctrl ]" (previously on row 8) slides down to row 9, and a
# brand-new "run selection" / "shift enter" entry takes over row 8.

# First, copy the bordered formatting that K8:L8 already has down onto
# K9:L9, so the relocated "debug cell" row keeps the boxed style instead of
# the plain (unbordered) style that row 9 currently carries.
$ws.Range("K8:L8").Copy()
$ws.Range("K9:L9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move "debug cell" / "ctrl ]" down into row 9.
$ws.Range("K9").Value = "debug cell"
$ws.Range("L9").Value = "ctrl ]"

# Put the new "run selection" / "shift enter" shortcut into the vacated K8:L8.
$ws.Range("K8").Value = "run selection"
$ws.Range("L8").Value = "shift enter"

# Widen columns K and L slightly so the newer, longer labels keep fitting
# (mirrors the bestFit column-width bump recorded in the saved file).
$ws.Columns.Item(11).ColumnWidth = 10.6
$ws.Columns.Item(12).ColumnWidth = 8.4

# The saved selection cursor moved from H21 to K21.
[void]$ws.Range("K21").Select()
